# Update to latest input spreadsheet
# - Re-express the JPY figures for one holding (row 31 of "ITR input data")
#   in USD (divide by the 107.92 JPY/USD rate) and correct its currency
#   code from the placeholder "megaJPY" to "USD".
# - Fix the region of another holding (row 37) from "Global" to "Asia".
# - Backfill missing market-cap (column J) figures for three holdings
#   (rows 45-47), formatted like the other "estimated" red figures
#   already used elsewhere in the sheet.
# - Restore the scroll position of the frozen panes on the "ITR input
#   data" and "ITR target input data" sheets.

$wb = $excel.ActiveWorkbook

$itrInput = $wb.Worksheets.Item("ITR input data")
$itrTarget = $wb.Worksheets.Item("ITR target input data")

# --- Row 31: currency was mislabeled "megaJPY" -- it's actually USD, and
# the financials for this row were entered in raw JPY instead of being
# converted to USD already (divide by the 107.92 JPY per USD rate used
# throughout the model).
$itrInput.Range("H31").Value = "USD"
$itrInput.Range("J31").Formula = "=879400*1000000/107.92"
$itrInput.Range("K31").Formula = "=5921500*1000000/107.92"
$itrInput.Range("N31").Formula = "=7444965*1000000/107.92"
$itrInput.Range("L31").Formula = "=M31+289549*1000000/107.92"
$itrInput.Range("M31").Formula = "=J31+2488741*1000000/107.92"

# --- Row 37: region should be "Asia", not "Global".
$itrInput.Range("E37").Value = "Asia"

# --- Rows 45-47: market_cap (column J) was missing; fill in the figures,
# matching the "estimated value" styling (red font) used for these rows.
$itrInput.Range("J45").Value = 50030000000
$itrInput.Range("J45").NumberFormat = "#,##0"
$itrInput.Range("J45").Font.Color = 255

$itrInput.Range("J46").Value = 590000000
$itrInput.Range("J46").Font.Color = 255

$itrInput.Range("J47").Value = 352130000
$itrInput.Range("J47").Font.Color = 255

# --- Restore view state: scroll the frozen panes back down to where the
# analyst was working, and update the active cell in "ITR input data".
$itrInput.Activate() | Out-Null
$itrInput.Range("J50").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 30
$excel.ActiveWindow.ScrollColumn = 2

$itrTarget.Activate() | Out-Null
$itrTarget.Range("L14").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 2

$itrInput.Activate() | Out-Null
